$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 299; existing rows 299:314 shift down to 300:315.
$ws.Rows.Item(299).Insert()

# Populate the newly inserted row 299 with a new weekly record (same fixed
# attributes as the other Betarraga / Vega Monumental Concepcion rows, with
# updated date and price figures).
$ws.Range("A299").Value = 11
$ws.Range("B299").Value = "Vega Monumental Concepción"
$ws.Range("C299").Value = "Bíobío"
$ws.Range("D299").Value = 44714
$ws.Range("D299").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E299").Value = 8
$ws.Range("F299").Value = 100114014
$ws.Range("G299").Value = "Betarraga"
$ws.Range("H299").Value = "Sin especificar"
$ws.Range("I299").Value = "Primera"
$ws.Range("J299").Value = 650
$ws.Range("K299").Value = 600
$ws.Range("L299").Value = 650
$ws.Range("M299").Value = 623
$ws.Range("N299").Value = "$/paquete 5 unidades"
$ws.Range("O299").Value = "Región Metropolitana"
$ws.Range("P299").Value = 125
$ws.Range("Q299").Value = 5
$ws.Range("R299").Value = "Hortaliza"
